# jobspersalary.xlsx update:
# Collapse the raw per-record rows (with a "Count" helper column) down to
# one summarized row per job title, with the salary column now holding the
# (weighted) average salary - computed with a formula - instead of a list
# of individual salary samples.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "Count" column's data (column C) - the column is no longer used.
$ws.Range("C1:C11").ClearContents()

# Drop the now-redundant detail rows 5-11; only the 3-job summary (rows 2-4)
# plus the header row remain.
$ws.Range("A5:B11").EntireRow.Delete()

# Header row
$ws.Range("A1").Value = "Jobs Name"
$ws.Range("B1").Value = "Salary (USD)"

# Data Scientist - single sample, unchanged
$ws.Range("A2").Value = "Data Scientist"
$ws.Range("B2").Value = 1350

# Data Engineer - average of the 3 salary samples (1350, 1350, 517)
$ws.Range("A3").Value = "Data Engineer"
$ws.Range("B3").Formula = "=(1350*2+517)/3"

# Data Analyst - average of the 5 salary samples (1350, 755, 755, 603.83, 1186)
$ws.Range("A4").Value = "Data Analyst"
$ws.Range("B4").Formula = "=(755*2+603.83+1186)/5"
